$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor data rows (timestamp, label is constant, ax, ay, az, gx, gy, gz)
# Each entry: "timestamp,ax,ay,az,gx,gy,gz"
$dataLines = @(
    "0,-12.57026290893555,-8.004177093505859,-9.469242095947266,-0.4692756063703969,-2.234323765550341,0.2831644466944987",
    "100,-11.83102798461914,-4.641905784606934,-9.06204319000244,0.5019868420703071,-2.09162780216762,-3.543231316975185",
    "200,5.961018562316895,-11.27237701416016,6.515813827514648,2.689777469634955,-1.050217630181979,-0.7848469325474046",
    "300,5.790210247039795,-16.61598014831543,5.962498664855957,-6.342504692077637,1.490384888648992,3.834645700454709",
    "400,-44.78383636474609,-8.077349662780762,-27.36338424682617,-5.673670666558469,3.692834442002409,2.466138475707613",
    "500,-42.66062545776367,-5.792407035827637,-30.21708297729492,3.314056737082347,6.166058949061821,0.5118652531078883",
    "600,-9.564473152160645,-21.24031639099121,-0.3582277297973633,2.383262787546406,-4.818590433469862,0.4682122468948429",
    "700,-10.97329330444336,-18.38039779663086,0.1557941436767578,8.378727106537124,0.4666124773877008,-1.115668596540174",
    "800,9.481355667114258,-9.28108024597168,9.836421966552734,-0.3396976134606717,0.4775580519012043,1.77553138222013",
    "900,10.69957828521728,-16.89757537841797,15.06509208679199,-5.64267150844849,3.583157830578955,3.719863629341134",
    "1000,-12.81566715240478,-8.839512825012207,-14.14643669128418,-1.766980731487296,16.75916714668267,-0.6191982626914752",
    "1100,-8.493080139160156,-3.757828950881958,-17.9176139831543,1.7845652954919,-3.40247355188642,4.17544686794281",
    "1200,-10.02292919158936,-2.351483345031738,5.807761192321777,8.421752619743332,-2.540225854941783,-3.182726718698209",
    "1300,-12.75368976593018,-6.235836982727051,6.349725723266602,0.1595347711017858,-0.04787222487585119,1.229899188450423",
    "1400,-79.57471466064453,-37.79908752441406,-28.62195587158203,-3.870756162064419,2.866504958271982,3.126545447962625",
    "1500,-83.42860412597656,-34.91292572021484,-30.07511520385743,-3.578196597099306,14.92317203113009,-0.7331702096121586",
    "1600,-3.41340184211731,-14.34111213684082,-5.178817749023438,0.1684449655669233,-0.0544045226915344,0.7549586423805885",
    "1700,-5.549047946929932,-10.40872192382812,-4.603366851806641,9.868156909942579,-0.5477856159210392,-1.925824952125532",
    "1800,13.1239423751831,-7.421818733215332,9.609323501586914,0.4101533753531434,-0.5063240064041949,1.204618217263898",
    "1900,13.5955753326416,-13.69983100891113,12.63335609436035,-4.432378809792649,1.422139526265,3.637724219049721",
    "2000,3.66697096824646,-4.846949577331543,-18.3593635559082,0.1461163333484263,9.942505402224434,-0.5480119875499323",
    "2100,8.207674980163574,-0.2672674059867859,-20.82013320922852,1.381819248199463,-3.363096782139367,0.9756205422537669",
    "2200,-3.444841146469116,0.0402781963348388,-6.625090599060059,5.626129477364711,-2.49930441464691,-2.115894814900024",
    "2300,-3.938364028930664,-3.24791407585144,-6.832444190979004,1.154986371312843,-0.3653818828718969,1.201646787779655",
    "2400,-71.62106323242188,-3.320010900497437,-19.81892585754395,-3.968483030796053,2.245044922828678,4.008097422122955",
    "2500,-75.16300201416016,-2.210659265518188,-18.6042366027832,-0.9093820844377891,8.189985391071852,-0.7311917849949312",
    "2600,-11.13384056091309,-14.18883228302002,-7.074286460876465,0.1449330214943292,1.165589673178594,-0.2754184092794358",
    "2700,-9.418439865112305,-9.676417350769045,-8.096508026123047,3.144946085555223,-6.092144768578668,-2.065853960173481",
    "2800,17.99904441833496,10.46211528778076,5.739192962646484,2.411846300108097,-2.998377110842354,-1.209750396864759",
    "2900,21.18538665771484,4.166534900665283,8.056502342224121,-1.357690974644261,0.8570090349231524,1.883484615598408"
)

$startRow = 2
for ($i = 0; $i -lt $dataLines.Count; $i++) {
    $fields = $dataLines[$i].Split(",")
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = [double]$fields[0]
    $ws.Cells.Item($row, 2).Value = "walkingToRunning"
    $ws.Cells.Item($row, 3).Value = [double]$fields[1]
    $ws.Cells.Item($row, 4).Value = [double]$fields[2]
    $ws.Cells.Item($row, 5).Value = [double]$fields[3]
    $ws.Cells.Item($row, 6).Value = [double]$fields[4]
    $ws.Cells.Item($row, 7).Value = [double]$fields[5]
    $ws.Cells.Item($row, 8).Value = [double]$fields[6]
}
